$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.341.49"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "3.182.49"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "594.18"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").Value = "148.75"
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.186.70"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").Value = "5.96"
$ws.Range("E11").Value = "  +6.15%  "
$ws.Range("D12").Value = "0.463"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "37.82"
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("D15").Value = "3.709.25"
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("D17").Value = "7.34"
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").Value = "3.180.60"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("D19").Value = "64.126.50"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "475.60"
$ws.Range("E20").Value = "  +3.40%  "
$ws.Range("D21").Value = "14.62"
$ws.Range("E21").Value = "  +3.04%  "
$ws.Range("D22").Value = "0.742"
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("D23").Value = "7.72"
$ws.Range("E23").Value = "  +4.65%  "
$ws.Range("E24").Value = "  +14.68%  "
$ws.Range("D25").Value = "13.26"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").Value = "81.69"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  +11.80%  "
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +3.89%  "
$ws.Range("D31").Value = "7.29"
$ws.Range("E31").Value = "  +5.99%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "0.118"
$ws.Range("E33").Value = "  +6.94%  "
$ws.Range("D34").Value = "28.46"
$ws.Range("E34").Value = "  +7.13%  "
$ws.Range("D35").Value = "0.0₃0863"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("E36").Value = "  +4.39%  "
$ws.Range("D37").Value = "6.25"
$ws.Range("E37").Value = "  +4.69%  "
$ws.Range("D38").Value = "2.33"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").Value = "475.25"
$ws.Range("E40").Value = "  +10.09%  "
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("D42").Value = "9.39"
$ws.Range("E42").Value = "  +8.12%  "
$ws.Range("D43").Value = "0.296"
$ws.Range("E43").Value = "  +9.95%  "
$ws.Range("E44").Value = "  +3.30%  "
$ws.Range("D45").Value = "2.923.49"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("D46").Value = "39.86"
$ws.Range("E46").Value = "  +11.03%  "
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("D48").Value = "133.20"
$ws.Range("E48").Value = "  +7.16%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "2.28"
$ws.Range("E50").Value = "  +7.11%  "
$ws.Range("E51").Value = "  +1.86%  "
